$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.158.86"
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("D3").Value = "2.266.55"
$ws.Range("E3").Value = "  +0.06%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'305.34"
$ws.Range("E5").Value = "  +0.49%  "
$ws.Range("D6").Value = "'96.11"
$ws.Range("E6").Value = "  +3.70%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +1.03%  "
$ws.Range("D10").Value = "'35.22"
$ws.Range("E10").Value = "  +8.60%  "
$ws.Range("D11").Value = "'0.0794"
$ws.Range("E11").Value = "  -0.36%  "
$ws.Range("E12").Value = "  -0.25%  "
$ws.Range("D13").Value = "'6.64"
$ws.Range("E13").Value = "  -0.25%  "
$ws.Range("D14").Value = "2.617.67"
$ws.Range("E14").Value = "  +0.07%  "
$ws.Range("E15").Value = "  +0.70%  "
$ws.Range("D16").Value = "2.269.40"
$ws.Range("E16").Value = "  -0.12%  "
$ws.Range("D17").Value = "'0.794"
$ws.Range("E17").Value = "  +1.25%  "
$ws.Range("D18").Value = "42.073.51"
$ws.Range("E18").Value = "  +0.77%  "
$ws.Range("E19").Value = "  -2.45%  "
$ws.Range("D20").Value = "0.0₃0908"
$ws.Range("E20").Value = "  +0.06%  "
$ws.Range("E21").Value = "  +0.80%  "
$ws.Range("D22").Value = "'67.76"
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("D23").Value = "'237.85"
$ws.Range("E23").Value = "  -2.60%  "
$ws.Range("E24").Value = "  -0.10%  "
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("D27").Value = "'23.71"
$ws.Range("E27").Value = "  -1.26%  "
$ws.Range("D28").Value = "'37.02"
$ws.Range("E28").Value = "  +5.79%  "
$ws.Range("D29").Value = "'9.52"
$ws.Range("E29").Value = "  -0.58%  "
$ws.Range("E30").Value = "  +2.12%  "
$ws.Range("D31").Value = "'159.41"
$ws.Range("E31").Value = "  +0.28%  "
$ws.Range("E32").Value = "  -0.57%  "
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("D34").Value = "'3.20"
$ws.Range("E34").Value = "  +5.89%  "
$ws.Range("D35").Value = "'0.0739"
$ws.Range("E35").Value = "  -0.59%  "
$ws.Range("D36").Value = "'17.06"
$ws.Range("E36").Value = "  +0.86%  "
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("E38").Value = "  -0.47%  "
$ws.Range("E39").Value = "  +2.12%  "
$ws.Range("E40").Value = "  -1.40%  "
$ws.Range("D41").Value = "'4.02"
$ws.Range("E41").Value = "  +2.76%  "
$ws.Range("E42").Value = "  +9.12%  "
$ws.Range("D43").Value = "1.988.59"
$ws.Range("E43").Value = "  -1.00%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'19.05"
$ws.Range("E44").Value = "  -4.54%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "'0.0284"
$ws.Range("E45").Value = "  +0.61%  "
$ws.Range("E46").Value = "  -5.03%  "
$ws.Range("E47").Value = "  +0.69%  "
$ws.Range("D48").Value = "'53.12"
$ws.Range("E48").Value = "  -0.05%  "
$ws.Range("E49").Value = "  +0.53%  "
$ws.Range("D50").Value = "'72.19"
$ws.Range("D51").Value = "'91.00"
$ws.Range("E51").Value = "  -0.76%  "
